$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: new SSD300 / voc0712 experiment on RTX 3070, batch_size = 4 ---
$ws.Range("A7").Value = "SSD300"
$ws.Range("B7").Value = "voc0712"
$ws.Range("C7").Value = "RTX 3070"
$ws.Range("D7").Value = "8G"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 4
$ws.Range("G7").Value = 24
$ws.Range("H7").Value = 8000
$ws.Range("M7").Value = "6D 9H"

# --- Row 8: new SSD300 / voc0712 experiment on RTX 3060 ---
$ws.Range("A8").Value = "SSD300"
$ws.Range("B8").Value = "voc0712"
$ws.Range("C8").Value = "RTX 3060"
$ws.Range("D8").Value = "12G"
$ws.Range("E8").Value = 2
$ws.Range("G8").Value = 24

# Match formatting of the existing data rows by copying each populated
# cell's format individually (avoids materialising blank styled cells in
# the untouched columns of rows 7/8).
$fmtCells = "A7","B7","C7","D7","E7","F7","G7","H7","M7","A8","B8","C8","D8","E8","F8","G8"
foreach ($addr in $fmtCells) {
    $ws.Range("A6").Copy()
    $ws.Range($addr).PasteSpecial(-4122)  # xlPasteFormats
}
$excel.CutCopyMode = $false

# Update the active selection to match the saved view state
$null = $ws.Range("J11").Select()
